$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "28.548.08"
Set-TextValue $ws.Range("E2") "  -1.70%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.962.62"
Set-TextValue $ws.Range("E3") "  -0.05%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.012"
Set-TextValue $ws.Range("E4") "  +0.68%  "

# Row 5
Set-TextValue $ws.Range("D5") "322.99"
Set-TextValue $ws.Range("E5") "  -1.21%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.010"
Set-TextValue $ws.Range("E6") "  +0.56%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4800"
Set-TextValue $ws.Range("E7") "  -3.90%  "

# Row 8
Set-TextValue $ws.Range("E8") "  -3.36%  "

# Row 9
Set-TextValue $ws.Range("D9") "54.25"
Set-TextValue $ws.Range("E9") "  +1.06%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.08492"
Set-TextValue $ws.Range("E10") "  -7.59%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -3.34%  "

# Row 12
Set-TextValue $ws.Range("D12") "22.43"
Set-TextValue $ws.Range("E12") "  -2.79%  "

# Row 13
Set-TextValue $ws.Range("B13") "Chainlink"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "7.592"
Set-TextValue $ws.Range("E13") "  -3.54%  "

# Row 14
Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "6.190"
Set-TextValue $ws.Range("E14") "  -3.85%  "

# Row 15
Set-TextValue $ws.Range("B15") "WrappedEther"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D15") "1.904.72"

# Row 16
Set-TextValue $ws.Range("E16") "  +0.77%  "

# Row 17
Set-TextValue $ws.Range("D17") "91.01"
Set-TextValue $ws.Range("E17") "  -0.44%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001073"
Set-TextValue $ws.Range("E18") "  -2.52%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06649"
Set-TextValue $ws.Range("E19") "  -0.34%  "

# Row 20
Set-TextValue $ws.Range("D20") "18.56"
Set-TextValue $ws.Range("E20") "  -3.15%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.011"
Set-TextValue $ws.Range("E21") "  +0.65%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.873"
Set-TextValue $ws.Range("E22") "  -1.16%  "

# Row 23
Set-TextValue $ws.Range("D23") "28.531.31"

# Row 24
Set-TextValue $ws.Range("E24") "  -4.29%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.301"
Set-TextValue $ws.Range("E25") "  +0.76%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.158.42"
Set-TextValue $ws.Range("E26") "  -2.36%  "

# Row 27
Set-TextValue $ws.Range("D27") "156.14"

# Row 28
Set-TextValue $ws.Range("E28") "  -1.25%  "

# Row 29
Set-TextValue $ws.Range("B29") "LidoDAOToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D29") "2.186"
Set-TextValue $ws.Range("E29") "  -3.51%  "

# Row 30
Set-TextValue $ws.Range("B30") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D30") "5.877"
Set-TextValue $ws.Range("E30") "  -4.91%  "

# Row 31
Set-TextValue $ws.Range("D31") "124.87"
Set-TextValue $ws.Range("E31") "  -1.57%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.9947"
Set-TextValue $ws.Range("E32") "  -4.59%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.09686"
Set-TextValue $ws.Range("E33") "  -1.74%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.464"
Set-TextValue $ws.Range("E34") "  -4.53%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.646"

# Row 36
Set-TextValue $ws.Range("D36") "3.696"
Set-TextValue $ws.Range("E36") "  +0.58%  "

# Row 37
Set-TextValue $ws.Range("D37") "9.116"
Set-TextValue $ws.Range("E37") "  +2.19%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.02341"
Set-TextValue $ws.Range("E38") "  -3.38%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.06255"
Set-TextValue $ws.Range("E39") "  -0.84%  "

# Row 40
Set-TextValue $ws.Range("E40") "  -2.93%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.6248"
Set-TextValue $ws.Range("E41") "  -3.28%  "

# Row 42
Set-TextValue $ws.Range("D42") "11.25"
Set-TextValue $ws.Range("E42") "  -1.64%  "

# Row 43
Set-TextValue $ws.Range("E43") "  +0.51%  "

# Row 44
Set-TextValue $ws.Range("E44") "  -3.68%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.356"
Set-TextValue $ws.Range("E45") "  +5.47%  "

# Row 46
Set-TextValue $ws.Range("B46") "EnergySwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "13.11"
Set-TextValue $ws.Range("E46") "  -2.23%  "

# Row 47
Set-TextValue $ws.Range("B47") "Decentraland"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.5968"
Set-TextValue $ws.Range("E47") "  -4.00%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -4.77%  "

# Row 49
Set-TextValue $ws.Range("D49") "3.413"
Set-TextValue $ws.Range("E49") "  -1.48%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.06845"
Set-TextValue $ws.Range("E50") "  -0.88%  "

# Row 51
Set-TextValue $ws.Range("E51") "  -1.05%  "
